$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.518534666666667
$ws.Range("N2").Value = 4.555604
$ws.Range("O2").Value = 0.1025715407499064
$ws.Range("P2").Value = 0.1025715407499064
$ws.Range("Q2").Value = 33.17963117898044
$ws.Range("R2").Value = 298.616680610824
$ws.Range("S2").Value = 0.005149110722311508
$ws.Range("T2").Value = 0.005149110722311508

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.418558666666667
$ws.Range("N3").Value = 13.255676
$ws.Range("O3").Value = 0.2984577041818288
$ws.Range("P3").Value = 0.2984577041818288
$ws.Range("Q3").Value = 96.54448470676178
$ws.Range("R3").Value = 868.900362360856
$ws.Range("S3").Value = 0.01498263313121319
$ws.Range("T3").Value = 0.01498263313121319

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.867545999999999
$ws.Range("N4").Value = 26.602638
$ws.Range("O4").Value = 0.5989707550682648
$ws.Range("P4").Value = 0.5989707550682649
$ws.Range("Q4").Value = 193.7538287410253
$ws.Range("R4").Value = 1743.784458669228
$ws.Range("S4").Value = 0.03006844505527073
$ws.Range("T4").Value = 0.03006844505527074

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.518534666666667
$ws.Range("N5").Value = 4.555604
$ws.Range("O5").Value = 0.1025715407499064
$ws.Range("P5").Value = 0.1025715407499064
$ws.Range("Q5").Value = 584.7155246868049
$ws.Range("R5").Value = 5262.439722181244
$ws.Range("S5").Value = 0.09074136362233495
$ws.Range("T5").Value = 0.09074136362233495

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.418558666666667
$ws.Range("N6").Value = 13.255676
$ws.Range("O6").Value = 0.2984577041818288
$ws.Range("P6").Value = 0.2984577041818288
$ws.Range("Q6").Value = 1701.37692991276
$ws.Range("R6").Value = 15312.39236921484
$ws.Range("S6").Value = 0.2640348274292187
$ws.Range("T6").Value = 0.2640348274292187

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.867545999999999
$ws.Range("N7").Value = 26.602638
$ws.Range("O7").Value = 0.5989707550682648
$ws.Range("P7").Value = 0.5989707550682649
$ws.Range("Q7").Value = 3414.470493094469
$ws.Range("R7").Value = 30730.23443785022
$ws.Range("S7").Value = 0.5298879463779875
$ws.Range("T7").Value = 0.5298879463779875

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.518534666666667
$ws.Range("N8").Value = 4.555604
$ws.Range("O8").Value = 0.1025715407499064
$ws.Range("P8").Value = 0.1025715407499064
$ws.Range("Q8").Value = 43.05118517810533
$ws.Range("R8").Value = 387.4606666029479
$ws.Range("S8").Value = 0.006681066405259907
$ws.Range("T8").Value = 0.006681066405259907

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.418558666666667
$ws.Range("N9").Value = 13.255676
$ws.Range("O9").Value = 0.2984577041818288
$ws.Range("P9").Value = 0.2984577041818288
$ws.Range("Q9").Value = 125.2682546896013
$ws.Range("R9").Value = 1127.414292206412
$ws.Range("S9").Value = 0.01944024362139686
$ws.Range("T9").Value = 0.01944024362139686

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.867545999999999
$ws.Range("N10").Value = 26.602638
$ws.Range("O10").Value = 0.5989707550682648
$ws.Range("P10").Value = 0.5989707550682649
$ws.Range("Q10").Value = 251.3991766545339
$ws.Range("R10").Value = 2262.592589890806
$ws.Range("S10").Value = 0.0390143636350066
$ws.Range("T10").Value = 0.03901436363500661
